# Apply textual replacements describing the diff between before.docx and after.docx.
# Each "old" string occurs exactly once in the document, so a simple
# Find/Replace (wdReplaceAll) for each pair is safe and deterministic.

$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-08 Sunday", "2024-12-09 Monday"),
    @("75×74=", "43×70="),
    @("67×90=", "37×34="),
    @("71×55=", "53×19="),
    @("56×36=", "81×79="),
    @("26×26=", "41×81="),
    @("65×31=", "33×20="),
    @("30×54=", "70×41="),
    @("19×46=", "56×83="),
    @("63×74=", "97×23="),
    @("83×79=", "20×22="),
    @("29×94=", "23×93="),
    @("21×83=", "90×74="),
    @("20×83=", "16×95="),
    @("68×44=", "91×87="),
    @("91×28=", "77×32="),
    @("70×42=", "18×66="),
    @("42×29=", "27×21="),
    @("38×79=", "92×68="),
    @("37×66=", "24×40="),
    @("80×15=", "17×85="),
    @("71×12=", "16×56="),
    @("48×13=", "99×82="),
    @("87×15=", "85×34="),
    @("51×63=", "24×37="),
    @("95×38=", "44×16=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
